$d = $word.ActiveDocument

# Common OOXML package wrapper used for targeted Range.InsertXML calls.
# Each call below replaces exactly one paragraph's content (identified by
# a Find match within it) with hand-built markup so the resulting runs /
# proofErr bookkeeping match the desired final shape precisely.

function Set-ParagraphXml($searchText, $paragraphXml) {
    $searchRange = $d.Content
    $found = $searchRange.Find.Execute($searchText)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $para = $searchRange.Paragraphs(1)
    $r = $para.Range
    $ooxml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $paragraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($ooxml) | Out-Null
}

# --- Edit 1: "Deployed URL shortener ..." -- collapse the 3 runs (split
# around the gramStart/gramEnd proofErr markers on "shortener") into a
# single run with the combined text. ---
$p1 = '<w:p w14:paraId="4D522C9E" w14:textId="7B04DEA8" w:rsidR="00911FDF" w:rsidRDefault="00911FDF" w:rsidP="006C03BD"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t>Deployed URL shortener for free and personal use amongst friends and family</w:t></w:r></w:p>'
Set-ParagraphXml "Deployed URL" $p1

# --- Edit 2: "Prime Access Consulting | Software Developer (...)" --
# update the end-date from "Present" to "March 2025", keeping the new
# month/year (and the closing paren) in their own italic runs. ---
$p2 = '<w:p w14:paraId="02EFD87A" w14:textId="0924FAFC" w:rsidR="00144061" w:rsidRPr="005F0687" w:rsidRDefault="00F70D46" w:rsidP="00532D90"><w:pPr><w:pStyle w:val="Heading3"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005F0687"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Prime Access Consulting | Software Developer (January 2023 -</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>March 2025</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>)</w:t></w:r></w:p>'
Set-ParagraphXml "Prime Access Consulting" $p2

# --- Edit 3: "Researched accessibility gaps ..." -- collapse the 3 runs
# (split around the gramStart/gramEnd proofErr markers on "alt") into a
# single run with the combined text. ---
$p3 = '<w:p w14:paraId="2B2017F0" w14:textId="711C78CA" w:rsidR="00126FBD" w:rsidRDefault="00126FBD" w:rsidP="00ED3D02"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:pPr><w:r w:rsidRPr="00126FBD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t>Researched accessibility gaps in the IIIF spec and common image viewers, drafting proposed metadata enhancements to support alt text and visual descriptions</w:t></w:r></w:p>'
Set-ParagraphXml "Researched accessibility" $p3

Write-Host "Done."
